# Update Malaysia economic data sheet with the latest reported figures.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: force a cell to stay text (matching the sheet's existing
# inline-string storage) before writing a numeric-looking string, so
# Excel doesn't silently coerce "4.47" -> 4.47 (number).
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
}

# Row 2 - Currency
Set-TextValue $ws.Range("C2") "4.47"
Set-TextValue $ws.Range("D2") "4.45"

# Row 3 - Stock Market
Set-TextValue $ws.Range("C3") "1610"
Set-TextValue $ws.Range("D3") "1612"

# Row 38 - Wages in Manufacturing
Set-TextValue $ws.Range("C38") "3442"
Set-TextValue $ws.Range("D38") "3414"
Set-TextValue $ws.Range("H38") "Sep/24"

# Row 61 - Foreign Exchange Reserves
Set-TextValue $ws.Range("C61") "117600"
Set-TextValue $ws.Range("D61") "119700"
Set-TextValue $ws.Range("H61") "Oct/24"

# Row 76 - Crude Oil Production
Set-TextValue $ws.Range("C76") "441"
Set-TextValue $ws.Range("D76") "462"
Set-TextValue $ws.Range("H76") "Jul/24"

# Row 113 - Construction Output
Set-TextValue $ws.Range("C113") "22.9"
Set-TextValue $ws.Range("D113") "20.2"
Set-TextValue $ws.Range("H113") "Sep/24"

Write-Host "Updates applied"
